# Weekly data refresh: a new price observation is inserted at row 104,
# shifting the existing rows 104-195 down to 105-196 (dimension grows to
# A1:T196). The new row carries this week's reading for the
# Tercera-quality, 16-unit-box Piña lot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104; Excel shifts rows 104:195 down to
# 105:196 and extends the used range / dimension automatically.
$ws.Rows(104).Insert()

$ws.Range("A104").Value = 4
$ws.Range("B104").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C104").Value = 'Los Lagos'
$ws.Range("D104").Value = 44586
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = 'Fruta'
$ws.Range("G104").Value = 100108
$ws.Range("H104").Value = 'Tropicales y subtropicales'
$ws.Range("I104").Value = 100108005
$ws.Range("J104").Value = 'Piña'
$ws.Range("K104").Value = 'Caramelo'
$ws.Range("L104").Value = 'Tercera'
$ws.Range("M104").Value = 200
$ws.Range("N104").Value = 19000
$ws.Range("O104").Value = 19000
$ws.Range("P104").Value = 19000
$ws.Range("Q104").Value = '$/caja 16 unidades'
$ws.Range("R104").Value = 'Ecuador'
$ws.Range("S104").Value = 1188
$ws.Range("T104").Value = 16
